$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two Paraguay - Primera Division rows (old rows 8 and 9).
# This shifts the rows below them (old 10, 11, 12) up to become rows 8, 9, 10.
$ws.Rows("8:9").Delete()

# Update odds that changed for the Argentina match (row 2).
$ws.Range("G2").Value = 2.8
$ws.Range("I2").Value = 3.2
$ws.Range("AI2").Value = 13
$ws.Range("AO2").Value = 21
$ws.Range("AR2").Value = 151
$ws.Range("AW2").Value = 4.75
$ws.Range("AX2").Value = 26

# Update odds that changed for the Colombia match (row 4).
$ws.Range("J4").Value = 2.88
$ws.Range("K4").Value = 1.95
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 5.5
$ws.Range("AA4").Value = 21
$ws.Range("AB4").Value = 41
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 17
$ws.Range("AM4").Value = 51

# Update odds that changed for the Mexico - Toluca match (row 6).
$ws.Range("G6").Value = 3
$ws.Range("I6").Value = 2.25
$ws.Range("J6").Value = 3.5
$ws.Range("L6").Value = 2.88
$ws.Range("X6").Value = 17
$ws.Range("AB6").Value = 29
$ws.Range("AD6").Value = 7
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 9
$ws.Range("AK6").Value = 21
$ws.Range("AX6").Value = 12

# Update odds for the Uruguay match, now at row 8 after the deletion.
$ws.Range("G8").Value = 4.5
$ws.Range("I8").Value = 1.85
$ws.Range("J8").Value = 4.75
$ws.Range("L8").Value = 2.6
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 3.25
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("W8").Value = 11
$ws.Range("AA8").Value = 41
$ws.Range("AC8").Value = 8
$ws.Range("AI8").Value = 8
$ws.Range("AK8").Value = 15
$ws.Range("AN8").Value = 6
$ws.Range("AR8").Value = 126
$ws.Range("AS8").Value = 301
$ws.Range("AV8").Value = 67
$ws.Range("AW8").Value = 3.75
$ws.Range("AX8").Value = 10
$ws.Range("AZ8").Value = 34
$ws.Range("BB8").Value = 201

# Update odds for the USA - MLS match, now at row 9 after the deletion.
$ws.Range("I9").Value = 7
$ws.Range("N9").Value = 26
$ws.Range("S9").Value = 1.18
$ws.Range("T9").Value = 4.5
$ws.Range("U9").Value = 1.53
$ws.Range("V9").Value = 2.38
$ws.Range("W9").Value = 13
$ws.Range("X9").Value = 10
$ws.Range("Y9").Value = 9.5
$ws.Range("AC9").Value = 26
$ws.Range("AG9").Value = 101
$ws.Range("AH9").Value = 29
$ws.Range("AN9").Value = 4
$ws.Range("AT9").Value = 4.5
$ws.Range("AU9").Value = 7.5
$ws.Range("AX9").Value = 29
$ws.Range("AY9").Value = 26
$ws.Range("BA9").Value = 81
$ws.Range("BB9").Value = 126
$ws.Range("BC9").Value = 251

# Update odds for the USA - USL Championship match, now at row 10 after the deletion.
$ws.Range("G10").Value = 1.62
$ws.Range("H10").Value = 4.05
$ws.Range("I10").Value = 4.3
$ws.Range("J10").Value = 2.12
$ws.Range("K10").Value = 2.45
$ws.Range("L10").Value = 4.4
$ws.Range("O10").Value = 1.17
$ws.Range("Q10").Value = 1.52
$ws.Range("R10").Value = 2.37
$ws.Range("S10").Value = 1.28
$ws.Range("T10").Value = 3.35
$ws.Range("X10").Value = 9.5
$ws.Range("AE10").Value = 13.5
$ws.Range("AI10").Value = 29
$ws.Range("AJ10").Value = 14.5
$ws.Range("AK10").Value = 70
$ws.Range("AL10").Value = 35
$ws.Range("AM10").Value = 32
$ws.Range("AO10").Value = 7.7
$ws.Range("AQ10").Value = 22
$ws.Range("AT10").Value = 3.35
$ws.Range("AU10").Value = 6.9
$ws.Range("AW10").Value = 6.5
$ws.Range("AX10").Value = 22
$ws.Range("AY10").Value = 24
$ws.Range("AZ10").Value = 110

